$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.0568700284933126
$ws.Range("E2").Value = 0.01211499930105764

$ws.Range("D3").Value = 0.02347435568743346
$ws.Range("E3").Value = 0.01145792177005145

$ws.Range("D4").Value = 0.03135589346508758
$ws.Range("E4").Value = 0.01369081574443798

$ws.Range("D5").Value = 0.03015879009749735
$ws.Range("E5").Value = 0.007095158597662632

$ws.Range("D6").Value = 0.03572628135440288
$ws.Range("E6").Value = 0.02069425901201605

$ws.Range("D7").Value = 0.01868395477373951
$ws.Range("E7").Value = 0.01390667659468847

$ws.Range("D8").Value = 0.004779801215917832
$ws.Range("E8").Value = 0.01850311850311837

$ws.Range("D9").Value = 0.006913652874472417
$ws.Range("E9").Value = 0.01092372556535071

$ws.Range("D10").Value = 0.07047474058895896
$ws.Range("E10").Value = -0.004512126339537459

$ws.Range("D11").Value = 0.07055423832227985
$ws.Range("E11").Value = -0.003943661971831047

$ws.Range("D12").Value = 0.1488091570789088
$ws.Range("E12").Value = -0.002350594771707293

$ws.Range("D13").Value = 0.3873964293593079
$ws.Range("E13").Value = -0.001308215593929796

$ws.Range("D14").Value = 0.1148026766886809
$ws.Range("E14").Value = 0.0004847309743092332

$ws.Range("E15").Value = 0.001367168893597004

$ws.Protect()
